# "running only one shareArticleLITest.java"
#
# The "Test Cases" sheet drives which JUnit tests actually run via a
# Runmode column (Y/N) and reflects the last-known Results (PASS/SKIP).
# This change flips every test back to "N"/"SKIP" except row 26
# (ShareArticleOnLITest), which is left as the only one enabled, and
# whose stale "SKIP" result is refreshed to "PASS". The active sheet
# selection/scroll position is also moved to track that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Runmode column (C): turn every test off ("N") apart from row 26,
# which stays "Y" (unchanged in the diff).
$ws.Range("C9").Value  = "N"   # CommentsMinMaxValidationTest
$ws.Range("C14").Value = "N"   # CommentsProfanityWordsCheckTest
$ws.Range("C15").Value = "N"   # UnsupportedTagsCommentsTest
$ws.Range("C20").Value = "N"   # ShareArticleOnTwitterTest
$ws.Range("C32").Value = "N"   # ShareArticleOnFBTest

# Results column (D): the Twitter-share test is no longer run, so its
# stale PASS is cleared back to SKIP, while the now-current LinkedIn
# test's stale SKIP is refreshed to PASS.
$ws.Range("D20").Value = "SKIP"  # ShareArticleOnTwitterTest
$ws.Range("D26").Value = "PASS"  # ShareArticleOnLITest

# Move the view / selection to the row that now matters (C26:C31,
# merged cell for ShareArticleOnLITest's Runmode), scrolled so row 9
# is at the top of the window.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C26:C31").Select()
